# Insert a new data row at row 27 (pushes the existing rows 27..153 down to 28..154),
# then populate the new row 27 with the new weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 27:153 down by one to make room for the new record.
$ws.Rows("27:27").Insert(-4121)

# Fill in the new row 27 - same as the record that used to occupy that row,
# except for an updated Fecha (Volumen) and Volumen value.
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = "Vega Modelo de Temuco"
$ws.Range("C27").Value = "La Araucanía"
$ws.Range("D27").Value = 45222
$ws.Range("E27").Value = 9
$ws.Range("F27").Value = 100112010
$ws.Range("G27").Value = "Achicoria"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = 10000
$ws.Range("N27").Value = "$/caja 18 unidades"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 556
$ws.Range("Q27").Value = 18
$ws.Range("R27").Value = "Hortaliza"
